# capstone_presentation.pptx -- "References" slide (slide 17):
#   - the references list (Content Placeholder 2) gains one more bulleted,
#     hyperlinked entry: "Coral Bleaching - Great Barrier Reef Foundation"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)

# Locate the "Content Placeholder 2" shape (the bulleted reference list).
$refShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $refShape = $candidate
    }
}

$tr = $refShape.TextFrame.TextRange
$beforeLen = $tr.Length

$newText = "Coral Bleaching - Great Barrier Reef Foundation"

# Typing a paragraph mark followed by the new text at the very end of the
# placeholder creates a new bulleted paragraph that inherits the preceding
# paragraph's bullet formatting (Wingdings "q" character bullet).
$tr.InsertAfter("`r" + $newText) | Out-Null

# Grab just the newly-typed run (skip the inserted paragraph mark) so we can
# give it its own hyperlink without touching the earlier entries.
$newRun = $tr.Characters($beforeLen + 2, $newText.Length)
$newRun.ActionSettings.Item(1).Hyperlink.Address = "https://www.barrierreef.org/the-reef/threats/coral-bleaching"
